$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 92
$ws.Range("I2").Value = 254
$ws.Range("J2").Value = 992
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 278
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 177
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 122
$ws.Range("T2").Value = 160
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 1609
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1590
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 13
